$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original (default) style so we can restore it after forcing text format
$origStyle = $ws.Cells.Item(2, 2).Style

# Column D price values: force text storage (some look like plain numbers)
# by setting the NumberFormat to Text before assignment, then restore the default style
# so no new/visible formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.959.11"
$ws.Range("D2").Style = $origStyle
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.463.01"
$ws.Range("D3").Style = $origStyle
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.02"
$ws.Range("D5").Style = $origStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.34"
$ws.Range("D6").Style = $origStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").Style = $origStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.461.77"
$ws.Range("D9").Style = $origStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.71"
$ws.Range("D11").Style = $origStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = $origStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.43"
$ws.Range("D14").Style = $origStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.902.54"
$ws.Range("D15").Style = $origStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.898.61"
$ws.Range("D16").Style = $origStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").Style = $origStyle
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.465.25"
$ws.Range("D18").Style = $origStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("D19").Style = $origStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.93"
$ws.Range("D20").Style = $origStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("D21").Style = $origStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.79"
$ws.Range("D22").Style = $origStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = $origStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.59"
$ws.Range("D24").Style = $origStyle
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.49"
$ws.Range("D27").Style = $origStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("D28").Style = $origStyle
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0803"
$ws.Range("D31").Style = $origStyle
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("D32").Style = $origStyle
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.36"
$ws.Range("D33").Style = $origStyle
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.397"
$ws.Range("D35").Style = $origStyle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.90"
$ws.Range("D36").Style = $origStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "366.69"
$ws.Range("D37").Style = $origStyle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.41"
$ws.Range("D38").Style = $origStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = $origStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.45"
$ws.Range("D42").Style = $origStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "150.17"
$ws.Range("D43").Style = $origStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.71"
$ws.Range("D44").Style = $origStyle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.52"
$ws.Range("D45").Style = $origStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("D46").Style = $origStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0958"
$ws.Range("D47").Style = $origStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("D48").Style = $origStyle
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0239"
$ws.Range("D49").Style = $origStyle
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0225"
$ws.Range("D50").Style = $origStyle
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.93"
$ws.Range("D51").Style = $origStyle

# Column E percentage-change text values (already non-numeric text, safe to set directly)
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("E3").Value = "  +5.74%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +3.29%  "
$ws.Range("E6").Value = "  +9.36%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("E9").Value = "  +5.74%  "
$ws.Range("E10").Value = "  +4.13%  "
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("E13").Value = "  +4.63%  "
$ws.Range("E14").Value = "  +11.56%  "
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("E16").Value = "  +4.23%  "
$ws.Range("E17").Value = "  +4.95%  "
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("E20").Value = "  +8.24%  "
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +8.03%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("E30").Value = "  +12.93%  "
$ws.Range("E31").Value = "  +9.41%  "
$ws.Range("E32").Value = "  +6.60%  "
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("E34").Value = "  +11.46%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").Value = "  +4.62%  "
$ws.Range("E37").Value = "  +13.70%  "
$ws.Range("E38").Value = "  +7.98%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("E41").Value = "  +11.04%  "
$ws.Range("E42").Value = "  +5.97%  "
$ws.Range("E43").Value = "  +8.57%  "
$ws.Range("E44").Value = "  +5.61%  "
$ws.Range("E45").Value = "  +5.77%  "
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("E49").Value = "  +9.18%  "
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("E51").Value = "  +5.22%  "
